$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5218.4
$ws.Range("I86").Value = 3364.6667
$ws.Range("K86").Value = 3364.6667
$ws.Range("M86").Value = -2241.6667
$ws.Range("H89").Value = 5218.4
$ws.Range("I89").Value = 3364.6667
$ws.Range("K89").Value = 16823.3335
$ws.Range("M89").Value = -11207.3335
$ws.Range("H103").Value = 1271
$ws.Range("J103").Value = 1249.5
$ws.Range("L103").Value = 3748.5
$ws.Range("N103").Value = -4920.5
$ws.Range("H125").Value = 1073.25
$ws.Range("I125").Value = 1055.1428
$ws.Range("J125").Value = 1200
$ws.Range("K125").Value = 9496.2852
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = -7036.2852
$ws.Range("N125").Value = -15720
$ws.Range("H127").Value = 1216.4375
$ws.Range("I127").Value = 1188.3636
$ws.Range("K127").Value = 3565.0908
$ws.Range("M127").Value = 1394.9092
$ws.Range("H132").Value = 2948.0833
$ws.Range("I132").Value = 2975.1714
$ws.Range("K132").Value = 8925.514200000001
$ws.Range("M132").Value = -6395.514200000001
$ws.Range("H135").Value = 6391.826
$ws.Range("I135").Value = 1474.3684
$ws.Range("K135").Value = 13269.3156
$ws.Range("M135").Value = -10734.3156
$ws.Range("H137").Value = 2498.36
$ws.Range("I137").Value = 2032.3334
$ws.Range("J137").Value = 3197.4
$ws.Range("K137").Value = 6097.0002
$ws.Range("L137").Value = 9592.200000000001
$ws.Range("M137").Value = -3547.0002
$ws.Range("N137").Value = -14692.2
$ws.Range("H138").Value = 2855.9768
$ws.Range("I138").Value = 1412.6428
$ws.Range("J138").Value = 3552.7585
$ws.Range("K138").Value = 4237.928400000001
$ws.Range("L138").Value = 10658.2755
$ws.Range("M138").Value = 902.0715999999993
$ws.Range("N138").Value = -20938.2755
$ws.Range("H141").Value = 2498.8
$ws.Range("I141").Value = 2607.7568
$ws.Range("J141").Value = 1155
$ws.Range("K141").Value = 7823.2704
$ws.Range("L141").Value = 3465
$ws.Range("M141").Value = -2643.2704
$ws.Range("N141").Value = -13825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2324.7727
$ws.Range("I74").Value = 2324.7727
$ws.Range("K74").Value = 2324.7727
$ws.Range("M74").Value = -1450.7727
$ws.Range("H77").Value = 2324.7727
$ws.Range("I77").Value = 2324.7727
$ws.Range("K77").Value = 11623.8635
$ws.Range("M77").Value = -7255.863499999999
$ws.Range("H92").Value = 35597.7
$ws.Range("J92").Value = 35597.7
$ws.Range("L92").Value = 35597.7
$ws.Range("N92").Value = -40589.7
$ws.Range("H95").Value = 20057.125
$ws.Range("J95").Value = 20057.125
$ws.Range("L95").Value = 20057.125
$ws.Range("N95").Value = -25549.125
$ws.Range("H113").Value = 19997.5
$ws.Range("J113").Value = 19997.5
$ws.Range("L113").Value = 19997.5
$ws.Range("N113").Value = -28675.5
$ws.Range("H132").Value = 2853.9524
$ws.Range("I132").Value = 2944.4285
$ws.Range("J132").Value = 2673
$ws.Range("K132").Value = 8833.2855
$ws.Range("L132").Value = 8019
$ws.Range("M132").Value = -6303.2855
$ws.Range("N132").Value = -13079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1653.9333
$ws.Range("I94").Value = 1488.7693
$ws.Range("K94").Value = 1488.7693
$ws.Range("M94").Value = -1037.7693
$ws.Range("H134").Value = 3039.5757
$ws.Range("I134").Value = 3076.2334
$ws.Range("K134").Value = 9228.700199999999
$ws.Range("M134").Value = -6693.700199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H58").Value = 2553.7368
$ws.Range("I58").Value = 2395.0625
$ws.Range("J58").Value = 3400
$ws.Range("K58").Value = 2395.0625
$ws.Range("L58").Value = 3400
$ws.Range("M58").Value = -2192.0625
$ws.Range("N58").Value = -3806
$ws.Range("H105").Value = 1758.8889
$ws.Range("I105").Value = 1390
$ws.Range("K105").Value = 1390
$ws.Range("M105").Value = 357
$ws.Range("H107").Value = 1489
$ws.Range("J107").Value = 1298.2
$ws.Range("L107").Value = 1298.2
$ws.Range("N107").Value = -5138.2
$ws.Range("H132").Value = 7854.2144
$ws.Range("I132").Value = 9380
$ws.Range("J132").Value = 5107.8
$ws.Range("K132").Value = 28140
$ws.Range("L132").Value = 15323.4
$ws.Range("M132").Value = -25610
$ws.Range("N132").Value = -20383.4
$ws.Range("H136").Value = 2553.7368
$ws.Range("I136").Value = 2395.0625
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 7185.1875
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -4635.1875
$ws.Range("N136").Value = -15300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14064
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40320
$ws.Range("H109").Value = 2682.2173
$ws.Range("I109").Value = 2359.55
$ws.Range("K109").Value = 7078.650000000001
$ws.Range("M109").Value = -6038.650000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 9999
$ws.Range("I35").Value = 9999
$ws.Range("K35").Value = 9999
$ws.Range("M35").Value = -9701
$ws.Range("H102").Value = 2999.6667
$ws.Range("I102").Value = 2999.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2999.6667
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = -1377.6667
$ws.Range("H122").Value = 1999
$ws.Range("J122").Value = 1999
$ws.Range("L122").Value = 5997
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 5339.6
$ws.Range("I126").Value = 2700
$ws.Range("J126").Value = 5999.5
$ws.Range("K126").Value = 8100
$ws.Range("L126").Value = 17998.5
$ws.Range("M126").Value = -5630
$ws.Range("N126").Value = -22938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 36028.855
$ws.Range("I7").Value = 38700.332
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 38700.332
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = -38588.332
$ws.Range("N7").Value = -20224
$ws.Range("H40").Value = 2999
$ws.Range("J40").Value = 2999
$ws.Range("L40").Value = 2999
$ws.Range("N40").Value = -3271
$ws.Range("H122").Value = 7455.231
$ws.Range("I122").Value = 11748.286
$ws.Range("J122").Value = 2446.6667
$ws.Range("K122").Value = 35244.858
$ws.Range("L122").Value = 7340.000100000001
$ws.Range("M122").Value = -32794.858
$ws.Range("N122").Value = -12240.0001
$ws.Range("H126").Value = 36028.855
$ws.Range("I126").Value = 38700.332
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 116100.996
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -113630.996
$ws.Range("N126").Value = -64940
$ws.Range("H132").Value = 2202.7446
$ws.Range("I132").Value = 1983.4324
$ws.Range("J132").Value = 3014.2
$ws.Range("K132").Value = 5950.2972
$ws.Range("L132").Value = 9042.599999999999
$ws.Range("M132").Value = -3420.2972
$ws.Range("N132").Value = -14102.6
$ws.Range("H136").Value = 1602.8889
$ws.Range("I136").Value = 1242.3572
$ws.Range("J136").Value = 2864.75
$ws.Range("K136").Value = 3727.0716
$ws.Range("L136").Value = 8594.25
$ws.Range("M136").Value = -1177.0716
$ws.Range("N136").Value = -13694.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3299.4443
$ws.Range("J122").Value = 3666
$ws.Range("L122").Value = 10998
$ws.Range("N122").Value = -15898
$ws.Range("H123").Value = 81999.5
$ws.Range("J123").Value = 81999.5
$ws.Range("L123").Value = 81999.5
$ws.Range("N123").Value = -91799.5
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H127").Value = 98997.5
$ws.Range("I127").Value = 98997.5
$ws.Range("K127").Value = 98997.5
$ws.Range("M127").Value = -94037.5
$ws.Range("H136").Value = 1313.8
$ws.Range("I136").Value = 1128.75
$ws.Range("J136").Value = 3287.6667
$ws.Range("K136").Value = 3386.25
$ws.Range("L136").Value = 9863.000100000001
$ws.Range("M136").Value = -836.25
$ws.Range("N136").Value = -14963.0001
